# Sync attendance_reports: swap order of "X, System" -> "System, X" in the
# "Recorded By" column (G) wherever the cell value is exactly
# "<email>, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "admin@admin.com, System") {
        $cell.Value = "System, admin@admin.com"
    }
}
